# Fruta / hortaliza, semanal
# Insert two new weekly-entry rows (Primera / Segunda) right after the
# existing row for 2021-04-16 (row 67 before the edit), pushing all the
# subsequent rows down by two and growing the sheet from 96 to 98 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 67 (this shifts the old row 67.. down to 69..)
$ws.Rows.Item(67).Insert()
$ws.Rows.Item(67).Insert()

# --- New row 67: "Primera" ---
$ws.Cells.Item(67,1).Value = 1
$ws.Cells.Item(67,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(67,3).Value = "Arica y Parinacota"
$ws.Cells.Item(67,4).Value2 = 44596
$ws.Cells.Item(67,5).Value = 15
$ws.Cells.Item(67,6).Value = 100112042
$ws.Cells.Item(67,7).Value = "Locoto"
$ws.Cells.Item(67,8).Value = "Sin especificar"
$ws.Cells.Item(67,9).Value = "Primera"
$ws.Cells.Item(67,10).Value = 120
$ws.Cells.Item(67,11).Value = 27000
$ws.Cells.Item(67,12).Value = 28000
$ws.Cells.Item(67,13).Value = 27500
$ws.Cells.Item(67,14).Value = "$/caja 20 kilos"
$ws.Cells.Item(67,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(67,16).Value = 1375
$ws.Cells.Item(67,17).Value = 20
$ws.Cells.Item(67,18).Value = "Hortaliza"

# --- New row 68: "Segunda" ---
$ws.Cells.Item(68,1).Value = 1
$ws.Cells.Item(68,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68,3).Value = "Arica y Parinacota"
$ws.Cells.Item(68,4).Value2 = 44596
$ws.Cells.Item(68,5).Value = 15
$ws.Cells.Item(68,6).Value = 100112042
$ws.Cells.Item(68,7).Value = "Locoto"
$ws.Cells.Item(68,8).Value = "Sin especificar"
$ws.Cells.Item(68,9).Value = "Segunda"
$ws.Cells.Item(68,10).Value = 200
$ws.Cells.Item(68,11).Value = 18000
$ws.Cells.Item(68,12).Value = 20000
$ws.Cells.Item(68,13).Value = 19000
$ws.Cells.Item(68,14).Value = "$/caja 20 kilos"
$ws.Cells.Item(68,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68,16).Value = 950
$ws.Cells.Item(68,17).Value = 20
$ws.Cells.Item(68,18).Value = "Hortaliza"
